# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
# Rows are matched by event name (column C) so the update is robust even if
# row positions differ between the two sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "南昌·第一届哥布林动漫游戏展——开学季&贺中秋" = 696
    "南昌·花绒万兽秋镜派对" = 32
    "上饶·星河城市动漫文化节" = 235
    "南昌·CM04动漫游戏博览会" = 2174
    "南昌·云芽动漫音乐嘉年华" = 3442
    "南昌·萌卡动漫展" = 858
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
